$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.531.48"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.603.22"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.06"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.71"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.52"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.062.23"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.428.62"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.629.42"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.53"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.11"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.31"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.22"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +6.20%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.848"
$ws.Range("E36").Value = "  +4.00%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.829"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.51"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.58"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.938.91"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.95"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  +1.57%  "
